$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 19300
$ws.Range("C2").Value = 225.2

$ws.Range("C2").Select()
